{"js": "// Apply the four textual edits described by the diff using search + replace.\nconst body = context.document.body;\n\nasync function replaceOnce(searchText, newText) {\n  const results = body.search(searchText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + searchText);\n  }\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 1) \"Esta caracter\u00edstica nos permite descomponer...\" -> \"...va a permitir descomponer...\"\nawait replaceOnce(\n  \"nos permite\",\n  \"va a permitir\"\n);\n\n// 2) \"ira construyendo\" -> \"construir\u00e1\"\nawait replaceOnce(\n  \"ira construyendo\",\n  \"construir\u00e1\"\n);\n\n// 3) \"Vamos a crear\" -> \"Se crear\u00e1\"\nawait replaceOnce(\n  \"Vamos a crear\",\n  \"Se crear\u00e1\"\n);\n\n// 4) \"vamos a ver c\u00f3mo se compone \" -> \"se visualizar\u00e1 en el navegador navegador c\u00f3mo se compone \"\nawait replaceOnce(\n  \"vamos a ver c\u00f3mo se compone \",\n  \"se visualizar\u00e1 en el navegador navegador c\u00f3mo se compone \"\n);\n", "ps1": "# Apply the four textual edits described by the diff using Find/Replace.\n$d = $word.ActiveDocument\n\nfunction Replace-Text($findText, $replaceText) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    # Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards, MatchSoundsLike,\n    #   MatchAllWordForms, Forward, Wrap, Format, ReplaceWith, Replace)\n    # MatchCase=$true keeps each replacement targeted at the exact phrase.\n    $find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n}\n\n# 1) \"Esta caracter\u00edstica nos permite descomponer...\" -> \"...va a permitir descomponer...\"\nReplace-Text \"nos permite\" \"va a permitir\"\n\n# 2) \"ira construyendo\" -> \"construir\u00e1\"\nReplace-Text \"ira construyendo\" \"construir\u00e1\"\n\n# 3) \"Vamos a crear\" -> \"Se crear\u00e1\"\nReplace-Text \"Vamos a crear\" \"Se crear\u00e1\"\n\n# 4) \"vamos a ver c\u00f3mo se compone \" -> \"se visualizar\u00e1 en el navegador navegador c\u00f3mo se compone \"\nReplace-Text \"vamos a ver c\u00f3mo se compone \" \"se visualizar\u00e1 en el navegador navegador c\u00f3mo se compone \"\n"}
